$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 2 (for the new falling/struggle/walkingToRunning
# gyroscope samples that were prepended to the dataset).
$ws.Rows.Item(2).Resize(3).Insert()

# Populate the newly inserted rows with the new sample data.
$newRows = @(
    @(-0.2070114476715818, -0.2780065764399136, 0.06705144135391006),
    @(-0.1988007093177122, -0.2540031636462492, 0.1641969842945828),
    @(-0.1353515688987338, -0.4769509890500239, -0.3330293473075421)
)

for ($i = 0; $i -lt $newRows.Length; $i++) {
    $r = 2 + $i
    $vals = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
}

# The last four rows of the old data (previously rows 19-22, now shifted to
# 22-25) are no longer part of the dataset, so remove them.
$ws.Rows.Item(22).Resize(4).Delete()
